$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct case dimensions for 01005 parts: update MPN/SPN to the correct
# (smaller/different case) part numbers.

# Row 3 (100n, C_0402 caps): GRM022R60J104KE15L -> GCM155R71C104JA55D
$ws.Range("F3").Value = "GCM155R71C104JA55D"
$ws.Range("H3").Value = "81-GCM155R71C104JA5D"

# Row 5 (10p, C4): GRM0225C1C100JA03L -> GRM1555C1E100GA01D
$ws.Range("F5").Value = "GRM1555C1E100GA01D"
$ws.Range("H5").Value = "81-GRM1555C1E100GA1D"

# Row 7 (1n, C14): strip stray trailing space from the Mouser SPN
$ws.Range("H7").Value = "81-GRM155R60J102KA1D"

# Restore the cursor/selection position left by the editing session
$ws.Range("H34").Select()
